# Update the three-digit ÷ one-digit division problems/answers in the table.
$d = $word.ActiveDocument

$replacements = @(
    @{ old = "870÷5=174, 0"; new = "924÷9=102, 6" },
    @{ old = "681÷9=75, 6";  new = "286÷4=71, 2" },
    @{ old = "731÷5=146, 1"; new = "991÷5=198, 1" },
    @{ old = "950÷2=475, 0"; new = "211÷7=30, 1" },
    @{ old = "805÷5=161, 0"; new = "929÷3=309, 2" },
    @{ old = "794÷5=158, 4"; new = "845÷9=93, 8" },
    @{ old = "754÷5=150, 4"; new = "259÷3=86, 1" },
    @{ old = "221÷6=36, 5";  new = "164÷3=54, 2" },
    @{ old = "576÷4=144, 0"; new = "324÷4=81, 0" },
    @{ old = "909÷5=181, 4"; new = "820÷4=205, 0" },
    @{ old = "501÷6=83, 3";  new = "867÷8=108, 3" },
    @{ old = "300÷3=100, 0"; new = "401÷5=80, 1" },
    @{ old = "528÷7=75, 3";  new = "817÷6=136, 1" },
    @{ old = "503÷8=62, 7";  new = "407÷8=50, 7" },
    @{ old = "994÷7=142, 0"; new = "458÷3=152, 2" },
    @{ old = "813÷4=203, 1"; new = "123÷9=13, 6" },
    @{ old = "487÷2=243, 1"; new = "406÷7=58, 0" },
    @{ old = "494÷6=82, 2";  new = "478÷6=79, 4" },
    @{ old = "833÷9=92, 5";  new = "843÷4=210, 3" },
    @{ old = "347÷5=69, 2";  new = "790÷2=395, 0" },
    @{ old = "687÷7=98, 1";  new = "547÷8=68, 3" },
    @{ old = "738÷2=369, 0"; new = "379÷9=42, 1" },
    @{ old = "962÷3=320, 2"; new = "278÷2=139, 0" },
    @{ old = "874÷7=124, 6"; new = "983÷5=196, 3" },
    @{ old = "696÷9=77, 3";  new = "820÷7=117, 1" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
